$d = $word.ActiveDocument

# Locate the last paragraph ("2.Create a program on to do list...") and
# collapse a range to its very end (just before the paragraph mark), then
# create a brand-new empty paragraph after it. Targeting InsertXML directly
# at the collapsed end-of-paragraph range would swallow the existing
# paragraph mark (and its text); inserting a fresh paragraph break first
# keeps "2.Create ..." intact and gives us an empty paragraph to fill in.
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($count)
$r = $lastPara.Range
$r.Collapse(0)
$r.InsertParagraphAfter()

$newCount = $d.Paragraphs.Count
$newPara = $d.Paragraphs.Item($newCount)
$r2 = $newPara.Range
$r2.Collapse(0)

# Insert the two new paragraphs (items 3 and 4) as raw OOXML so we can
# reproduce the exact run-splitting / proofErr markers Word's live
# grammar & spell checker would have produced while the author typed.
$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:rPr>
      <w:sz w:val="36"/>
      <w:szCs w:val="36"/>
    </w:rPr>
  </w:pPr>
  <w:proofErr w:type="gramStart"/>
  <w:r>
    <w:rPr>
      <w:sz w:val="36"/>
      <w:szCs w:val="36"/>
    </w:rPr>
    <w:t>3.Create</w:t>
  </w:r>
  <w:proofErr w:type="gramEnd"/>
  <w:r>
    <w:rPr>
      <w:sz w:val="36"/>
      <w:szCs w:val="36"/>
    </w:rPr>
    <w:t xml:space="preserve"> a program on digital clock</w:t>
  </w:r>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:rPr>
      <w:sz w:val="36"/>
      <w:szCs w:val="36"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:sz w:val="36"/>
      <w:szCs w:val="36"/>
    </w:rPr>
    <w:t xml:space="preserve">4.Create a program on </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:sz w:val="36"/>
      <w:szCs w:val="36"/>
    </w:rPr>
    <w:t>count down</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr>
      <w:sz w:val="36"/>
      <w:szCs w:val="36"/>
    </w:rPr>
    <w:t>(birthday/new year)</w:t>
  </w:r>
</w:p>
'@

$r2.InsertXML($xml)
